$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. _set_SCENARIOS - content unchanged, just move the selection to B2
# ---------------------------------------------------------------------------
$wsScenarios = $wb.Worksheets.Item("_set_SCENARIOS")
$wsScenarios.Range("B2").Select()

# ---------------------------------------------------------------------------
# 2. _set_TECHS - collapse the technology list down to a single "Power plant"
# ---------------------------------------------------------------------------
$wsTechs = $wb.Worksheets.Item("_set_TECHS")
$wsTechs.Range("A2").Value = "Power plant"
$wsTechs.Range("A3:A4").ClearContents()
$wsTechs.Range("C11").Select()

# ---------------------------------------------------------------------------
# 3. _set_FLOWS - drop the flows_Agg column and collapse to a single
#    "Electricity" flow
# ---------------------------------------------------------------------------
$wsFlows = $wb.Worksheets.Item("_set_FLOWS")
$wsFlows.Range("A2").Value = "Electricity"
$wsFlows.Range("A3:A4").ClearContents()
$wsFlows.Columns.Item(2).Delete()
$wsFlows.Range("E10").Select()

# ---------------------------------------------------------------------------
# 4. _set_FLOWS_AGG - no longer needed, remove the sheet entirely
# ---------------------------------------------------------------------------
$wsFlowsAgg = $wb.Worksheets.Item("_set_FLOWS_AGG")
$wsFlowsAgg.Delete()

# ---------------------------------------------------------------------------
# 5. _set_YEARS - content unchanged, just move the selection
# ---------------------------------------------------------------------------
$wsYears = $wb.Worksheets.Item("_set_YEARS")
$wsYears.Range("G20:G21").Select()

# ---------------------------------------------------------------------------
# 6. _set_LOADFACTORS - content unchanged, just move the selection
# ---------------------------------------------------------------------------
$wsLoadfactors = $wb.Worksheets.Item("_set_LOADFACTORS")
$wsLoadfactors.Range("C3").Select()

# ---------------------------------------------------------------------------
# 7. _set_Costs - rename the cost items and split "o&m" into its own labelled
#    "O&M costs" row (investment -> Investment costs, o&m -> O&M costs)
# ---------------------------------------------------------------------------
$wsCosts = $wb.Worksheets.Item("_set_Costs")
$wsCosts.Range("A2").Value = "Investment costs"
$wsCosts.Range("B2").Value = "inv"
$wsCosts.Range("A3").Value = "O&M costs"
$wsCosts.Range("B3").Value = "om"
$wsCosts.Activate()
$wsCosts.Range("B4").Select()
